# Weekly update: insert a new daily record for "Vega Modelo de Temuco - Mango"
# above the current row 179, pushing all subsequent rows down by one.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 179; rows 179..199 shift down to 180..200.
$ws.Rows.Item(179).Insert()

# Populate the newly inserted row 179 with the new daily data point.
# Columns that are constant across every record in this sheet are copied
# from the row immediately below (now row 180, the former row 179).
$ws.Cells.Item(179, 1).Value = 10                               # A: Mercado ID
$ws.Cells.Item(179, 2).Value = "Vega Modelo de Temuco"          # B: Mercado
$ws.Cells.Item(179, 3).Value = "La Araucanía"                   # C: Región
$ws.Cells.Item(179, 4).Value = 44449                            # D: Fecha
$ws.Cells.Item(179, 5).Value = 9                                # E: Codreg
$ws.Cells.Item(179, 6).Value = "Fruta"                          # F: Tipo
$ws.Cells.Item(179, 7).Value = 100108                           # G: Producto ID
$ws.Cells.Item(179, 8).Value = "Tropicales y subtropicales"     # H: Producto
$ws.Cells.Item(179, 9).Value = 100108002                        # I: Categoría ID
$ws.Cells.Item(179, 10).Value = "Mango"                         # J: Categoría
$ws.Cells.Item(179, 11).Value = "Sin especificar"                # K: Variedad
$ws.Cells.Item(179, 12).Value = "Primera"                        # L: Calidad
$ws.Cells.Item(179, 13).Value = 500                              # M: Volumen
$ws.Cells.Item(179, 14).Value = 9000                             # N: Precio mínimo
$ws.Cells.Item(179, 15).Value = 9000                             # O: Precio máximo
$ws.Cells.Item(179, 16).Value = 9000                             # P: Precio promedio ponderado
$ws.Cells.Item(179, 17).Value = "$/bandeja 4 kilos"              # Q: Unidad de comercialización
$ws.Cells.Item(179, 18).Value = "Brasil"                         # R: Origen
$ws.Cells.Item(179, 19).Value = 2250                             # S: Precio $/Kg
$ws.Cells.Item(179, 20).Value = 4                                # T: Kg / unidad
